$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 already has A5 (Sno=4) and B5 (IP="120.221.212.160") populated from a
# previous edit. Fill in the remaining columns for this new VirusTotal record.
$ws.Range("C5").Value = "Malicious"

# D5: hyperlink cell, mirrors D2:D4 (Hyperlink style + link to the VirusTotal
# detection page for the new IP).
$linkText = "https://www.virustotal.com/gui/ip-address/120.221.212.160/detection"
$ws.Range("D5").Value = $linkText
$ws.Hyperlinks.Add($ws.Range("D5"), $linkText) | Out-Null
$ws.Range("D5").Style = "Hyperlink"

$ws.Range("E5").Value = "{'harmless': 56, 'malicious': 12, 'suspicious': 1, 'undetected': 21, 'timeout': 0}"
$ws.Range("F5").Value = "China"

# G5 / H5: date/time values, same numeric (serial) storage + number format as G2:H4
$ws.Range("G5").Value = 45328.57457175926
$ws.Range("G5").NumberFormat = $ws.Range("G4").NumberFormat
$ws.Range("H5").Value = 45340.69390046296
$ws.Range("H5").NumberFormat = $ws.Range("H4").NumberFormat

$ws.Range("I5").Value = "China Mobile Communications Group Co., Ltd."
